# Apply the "updated the excel reformatting and dynamic recipient" edit.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("bookingData")
$ws2 = $wb.Worksheets.Item("giftData")

# --- giftData sheet: add quantity column, change values, add hyperlink ---

# Row 2 data updates

# A2: reformat from Text to General, then overwrite with the new numeric amount
$ws2.Range("A2").NumberFormat = "General"
$ws2.Range("A2").Value = 2000

# C2: mobile number becomes a text value (was a raw number before)
$ws2.Range("C2").NumberFormat = "@"
$ws2.Range("C2").Value = "9900443322"

# Header row: new column E = "quantity"  (added to shared strings after the mobile number)
$ws2.Range("E1").Value = "quantity"

# E2: quantity value, stored as text
$ws2.Range("E2").NumberFormat = "@"
$ws2.Range("E2").Value = "5"

# Turn D2 (sssdd) into a mailto hyperlink, displaying the same text
$ws2.Hyperlinks.Add($ws2.Range("D2"), "mailto:sssdd@gmail.com", [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "sssdd")

# --- bookingData sheet: just a selection/view change ---
$ws1.Select()
$ws1.Range("F3").Select()

# giftData stays the active/tab-selected sheet, with the cursor on the new cell
$ws2.Select()
$ws2.Range("E2").Select()

$wb.Save()
